# Daily "remaining days" refresh for the tracking sheet.
#
# Columns: D = total days (总天), E = remaining days (剩余), F = cycle
# start date (开始时间) stored as a plain yyyyMMdd integer.
#
# For every data row (2..99):
#   - skip rows that are already at a full/fresh cycle (remaining == total)
#   - if remaining has hit 1, the cycle is over: reset remaining to the
#     total, and roll the start date forward by `total` days (today)
#   - otherwise just tick remaining down by one day

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 99

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $totalDays = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $startDate = $ws.Cells.Item($r, 6).Value2

    if ($null -eq $remaining) {
        continue
    }

    if ($remaining -eq $totalDays) {
        # Already a fresh/full cycle - nothing to do today.
        continue
    }

    if ($remaining -eq 1) {
        # Cycle expired: start a new one today and reset the countdown.
        $cycleStart = [datetime]::ParseExact([string]$startDate, "yyyyMMdd", $null)
        $newCycleStart = $cycleStart.AddDays([double]$totalDays)

        $ws.Cells.Item($r, 5).Value = $totalDays
        $ws.Cells.Item($r, 6).Value = [int]$newCycleStart.ToString("yyyyMMdd")
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
